# Auto update Excel log 2026-02-04 14:34:50
# Appends new sensor log rows to the PIR, Humidity, and Temperature sheets,
# mirroring the continuously-growing sensor log export.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: append rows 393-404
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PIR")

$pirData = @(
    @("393","14:33:48","14:00","Bathroom","No Motion","Inactive"),
    @("394","14:33:51","14:00","Bathroom","No Motion","Inactive"),
    @("395","14:33:52","14:00","Bathroom","Motion Detected","Active"),
    @("396","14:33:58","14:00","Bathroom","No Motion","Inactive"),
    @("397","14:34:03","14:00","Bathroom","No Motion","Inactive"),
    @("398","14:34:08","14:00","Bathroom","No Motion","Inactive"),
    @("399","14:34:13","14:00","Bathroom","No Motion","Inactive"),
    @("400","14:34:16","14:00","Bathroom","Motion Detected","Active"),
    @("401","14:34:25","14:00","Bathroom","No Motion","Inactive"),
    @("402","14:34:30","14:00","Bathroom","No Motion","Inactive"),
    @("403","14:34:31","14:00","Bathroom","Motion Detected","Active"),
    @("404","14:34:41","14:00","Bathroom","No Motion","Inactive")
)

foreach ($rec in $pirData) {
    $row = $rec[0]
    # Leading apostrophe forces Excel to keep the date-looking value as literal
    # text instead of auto-converting it to a date serial number.
    $ws.Cells.Item($row, 1).Value = "'2026-02-04"
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
}

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 320-332
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Humidity")

$humidityData = @(
    @("320","14:33:44","79.8%"),
    @("321","14:33:46","79.0%"),
    @("322","14:33:49","79.9%"),
    @("323","14:33:53","79.1%"),
    @("324","14:33:58","80.0%"),
    @("325","14:34:03","79.0%"),
    @("326","14:34:08","79.9%"),
    @("327","14:34:18","80.0%"),
    @("328","14:34:23","79.2%"),
    @("329","14:34:28","80.2%"),
    @("330","14:34:33","79.2%"),
    @("331","14:34:38","80.0%"),
    @("332","14:34:43","79.3%")
)

foreach ($rec in $humidityData) {
    $row = $rec[0]
    $ws.Cells.Item($row, 1).Value = "'2026-02-04"
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = "14:00"
    $ws.Cells.Item($row, 4).Value = "Bathroom"
    # Leading apostrophe keeps the percentage-looking value as literal text
    # instead of Excel auto-converting it into a numeric percentage.
    $ws.Cells.Item($row, 5).Value = "'" + $rec[2]
    $ws.Cells.Item($row, 6).Value = "Active"
}

# ---------------------------------------------------------------------------
# Temperature sheet: append rows 320-332
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Temperature")

$temperatureData = @(
    @("320","14:33:45","24.3C"),
    @("321","14:33:47","24.3C"),
    @("322","14:33:50","24.2C"),
    @("323","14:33:54","24.2C"),
    @("324","14:34:00","24.2C"),
    @("325","14:34:04","24.2C"),
    @("326","14:34:10","24.2C"),
    @("327","14:34:19","24.2C"),
    @("328","14:34:24","24.2C"),
    @("329","14:34:29","24.2C"),
    @("330","14:34:34","24.2C"),
    @("331","14:34:39","24.2C"),
    @("332","14:34:44","24.2C")
)

foreach ($rec in $temperatureData) {
    $row = $rec[0]
    $ws.Cells.Item($row, 1).Value = "'2026-02-04"
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = "14:00"
    $ws.Cells.Item($row, 4).Value = "Bathroom"
    $ws.Cells.Item($row, 5).Value = $rec[2]
    $ws.Cells.Item($row, 6).Value = "Active"
}
